# Apply weekly rotation of Fecha/Volumen/Precio columns (D, J, K, L, M, P)
# across rows 2, 7, 13, 8, 5, 11, 4 (in that cyclic order).
# new_row[i] = old_row[i+1] (wrapping), i.e.:
#   new 2 = old 7, new 7 = old 13, new 13 = old 8, new 8 = old 5,
#   new 5 = old 11, new 11 = old 4, new 4 = old 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "P")
$cycle = @(2, 7, 13, 8, 5, 11, 4)

# Snapshot the original values for each row/column involved before any writes.
$orig = @{}
foreach ($r in $cycle) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Apply the rotation: each row takes the values from the next row in the cycle.
$n = $cycle.Length
for ($i = 0; $i -lt $n; $i++) {
    $destRow = $cycle[$i]
    $srcRow = $cycle[($i + 1) % $n]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $orig[$srcRow][$col]
    }
}
